$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 4: mention additional frameworks used in the course.
$tr.Paragraphs(4).Text = "We will use Python, Keras, Tensorflow, HuggingFace Encoder/Decoder, Google PaLM (Pathway Language Model), Llama (Large Large Model Meta AI), and PyTorch (Optional)."

# Paragraph 6: merge the two runs into a single continuous sentence.
$tr.Paragraphs(6).Text = "1. Concepts in mathematics and statistics: In Machine Learning class, we use numpy library and scikit-learn library for matrix operation and statistical regression modeling. In Deep Learning, we will use them to convert the dataset into tensor."

# Paragraph 8: fix typo "porojects" -> "projects".
$tr.Paragraphs(8).Text = "For projects:"
